$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$ws.Range("B2").Value = 0.02301637451559633
$ws.Range("C2").Value = 1.453483806356359
$ws.Range("D2").Value = 6.459978470793084
$ws.Range("E2").Value = 2.541648770147654
$ws.Range("F2").Value = 2.604308190894869
